$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.813278666666667
$ws.Range("H2").Value = 26.439836
$ws.Range("I2").Value = 0.3770976991891536
$ws.Range("J2").Value = 0.3770976991891536
$ws.Range("M2").Value = 10.46510533333333
$ws.Range("N2").Value = 31.395316
$ws.Range("O2").Value = 0.5554075997074865
$ws.Range("P2").Value = 0.5554075997074865
$ws.Range("Q2").Value = 92.23188957868624
$ws.Range("R2").Value = 830.087006208176
$ws.Range("S2").Value = 0.2094429279618636
$ws.Range("T2").Value = 0.2094429279618636

# Row 3
$ws.Range("G3").Value = 8.813278666666667
$ws.Range("H3").Value = 26.439836
$ws.Range("I3").Value = 0.3770976991891536
$ws.Range("J3").Value = 0.3770976991891536
$ws.Range("O3").Value = 0.3053945925621632
$ws.Range("P3").Value = 0.3053945925621632
$ws.Range("Q3").Value = 50.7143228755889
$ws.Range("R3").Value = 456.4289058803001
$ws.Range("S3").Value = 0.1151635982000007
$ws.Range("T3").Value = 0.1151635982000007

# Row 4
$ws.Range("G4").Value = 8.813278666666667
$ws.Range("H4").Value = 26.439836
$ws.Range("I4").Value = 0.3770976991891536
$ws.Range("J4").Value = 0.3770976991891536
$ws.Range("M4").Value = 2.146766
$ws.Range("N4").Value = 6.440298
$ws.Range("O4").Value = 0.1139338891693565
$ws.Range("P4").Value = 0.1139338891693565
$ws.Range("Q4").Value = 18.92004699012534
$ws.Range("R4").Value = 170.280422911128
$ws.Range("S4").Value = 0.04296420746543637
$ws.Range("T4").Value = 0.04296420746543638

# Row 5
$ws.Range("G5").Value = 8.813278666666667
$ws.Range("H5").Value = 26.439836
$ws.Range("I5").Value = 0.3770976991891536
$ws.Range("J5").Value = 0.3770976991891536
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4760280000000001
$ws.Range("N5").Value = 1.428084
$ws.Range("O5").Value = 0.02526391856099382
$ws.Range("P5").Value = 0.02526391856099382
$ws.Range("Q5").Value = 4.195367417136
$ws.Range("R5").Value = 37.758306754224
$ws.Range("S5").Value = 0.009526965561852924
$ws.Range("T5").Value = 0.009526965561852924

# Row 6
$ws.Range("I6").Value = 0.5522024902836482
$ws.Range("J6").Value = 0.5522024902836482
$ws.Range("M6").Value = 10.46510533333333
$ws.Range("N6").Value = 31.395316
$ws.Range("O6").Value = 0.5554075997074865
$ws.Range("P6").Value = 0.5554075997074865
$ws.Range("Q6").Value = 135.0596389700325
$ws.Range("R6").Value = 1215.536750730292
$ws.Range("S6").Value = 0.3066974596809377
$ws.Range("T6").Value = 0.3066974596809377

# Row 7
$ws.Range("I7").Value = 0.5522024902836482
$ws.Range("J7").Value = 0.5522024902836482
$ws.Range("O7").Value = 0.3053945925621632
$ws.Range("P7").Value = 0.3053945925621632
$ws.Range("S7").Value = 0.1686396545319866
$ws.Range("T7").Value = 0.1686396545319866

# Row 8
$ws.Range("I8").Value = 0.5522024902836482
$ws.Range("J8").Value = 0.5522024902836482
$ws.Range("M8").Value = 2.146766
$ws.Range("N8").Value = 6.440298
$ws.Range("O8").Value = 0.1139338891693565
$ws.Range("P8").Value = 0.1139338891693565
$ws.Range("Q8").Value = 27.70554444298067
$ws.Range("R8").Value = 249.349899986826
$ws.Range("S8").Value = 0.06291457732701985
$ws.Range("T8").Value = 0.06291457732701985

# Row 9
$ws.Range("I9").Value = 0.5522024902836482
$ws.Range("J9").Value = 0.5522024902836482
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4760280000000001
$ws.Range("N9").Value = 1.428084
$ws.Range("O9").Value = 0.02526391856099382
$ws.Range("P9").Value = 0.02526391856099382
$ws.Range("Q9").Value = 6.143480430612001
$ws.Range("R9").Value = 55.29132387550801
$ws.Range("S9").Value = 0.01395079874370407
$ws.Range("T9").Value = 0.01395079874370407

# Row 10
$ws.Range("G10").Value = 1.649921333333333
$ws.Range("H10").Value = 4.949764
$ws.Range("I10").Value = 0.07059592260441032
$ws.Range("J10").Value = 0.07059592260441033
$ws.Range("M10").Value = 10.46510533333333
$ws.Range("N10").Value = 31.395316
$ws.Range("O10").Value = 0.5554075997074865
$ws.Range("P10").Value = 0.5554075997074865
$ws.Range("Q10").Value = 17.26660054504711
$ws.Range("R10").Value = 155.399404905424
$ws.Range("S10").Value = 0.03920951192285103
$ws.Range("T10").Value = 0.03920951192285103

# Row 11
$ws.Range("G11").Value = 1.649921333333333
$ws.Range("H11").Value = 4.949764
$ws.Range("I11").Value = 0.07059592260441032
$ws.Range("J11").Value = 0.07059592260441033
$ws.Range("O11").Value = 0.3053945925621632
$ws.Range("P11").Value = 0.3053945925621632
$ws.Range("Q11").Value = 9.494156077744446
$ws.Range("R11").Value = 85.44740469970002
$ws.Range("S11").Value = 0.02155961302032389
$ws.Range("T11").Value = 0.0215596130203239

# Row 12
$ws.Range("G12").Value = 1.649921333333333
$ws.Range("H12").Value = 4.949764
$ws.Range("I12").Value = 0.07059592260441032
$ws.Range("J12").Value = 0.07059592260441033
$ws.Range("M12").Value = 2.146766
$ws.Range("N12").Value = 6.440298
$ws.Range("O12").Value = 0.1139338891693565
$ws.Range("P12").Value = 0.1139338891693565
$ws.Range("Q12").Value = 3.541995021074666
$ws.Range("R12").Value = 31.877955189672
$ws.Range("S12").Value = 0.008043268021819356
$ws.Range("T12").Value = 0.008043268021819358

# Row 13
$ws.Range("G13").Value = 1.649921333333333
$ws.Range("H13").Value = 4.949764
$ws.Range("I13").Value = 0.07059592260441032
$ws.Range("J13").Value = 0.07059592260441033
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4760280000000001
$ws.Range("N13").Value = 1.428084
$ws.Range("O13").Value = 0.02526391856099382
$ws.Range("P13").Value = 0.02526391856099382
$ws.Range("Q13").Value = 0.7854087524640001
$ws.Range("R13").Value = 7.068678772176001
$ws.Range("S13").Value = 0.001783529639416045
$ws.Range("T13").Value = 0.001783529639416045

# Row 14
$ws.Range("G14").Value = 0.002428
$ws.Range("H14").Value = 0.007284
$ws.Range("I14").Value = 0.0001038879227879399
$ws.Range("J14").Value = 0.0001038879227879399
$ws.Range("M14").Value = 10.46510533333333
$ws.Range("N14").Value = 31.395316
$ws.Range("O14").Value = 0.5554075997074865
$ws.Range("P14").Value = 0.5554075997074865
$ws.Range("Q14").Value = 0.02540927574933333
$ws.Range("R14").Value = 0.228683481744
$ws.Range("S14").Value = 0.00005770014183424642
$ws.Range("T14").Value = 0.00005770014183424642

# Row 15
$ws.Range("G15").Value = 0.002428
$ws.Range("H15").Value = 0.007284
$ws.Range("I15").Value = 0.0001038879227879399
$ws.Range("J15").Value = 0.0001038879227879399
$ws.Range("O15").Value = 0.3053945925621632
$ws.Range("P15").Value = 0.3053945925621632
$ws.Range("Q15").Value = 0.01397146063333334
$ws.Range("R15").Value = 0.1257431457
$ws.Range("S15").Value = 0.00003172680985195239
$ws.Range("T15").Value = 0.00003172680985195239

# Row 16
$ws.Range("G16").Value = 0.002428
$ws.Range("H16").Value = 0.007284
$ws.Range("I16").Value = 0.0001038879227879399
$ws.Range("J16").Value = 0.0001038879227879399
$ws.Range("M16").Value = 2.146766
$ws.Range("N16").Value = 6.440298
$ws.Range("O16").Value = 0.1139338891693565
$ws.Range("P16").Value = 0.1139338891693565
$ws.Range("Q16").Value = 0.005212347848
$ws.Range("R16").Value = 0.046911130632
$ws.Range("S16").Value = 0.00001183635508095582
$ws.Range("T16").Value = 0.00001183635508095582

# Row 17
$ws.Range("G17").Value = 0.002428
$ws.Range("H17").Value = 0.007284
$ws.Range("I17").Value = 0.0001038879227879399
$ws.Range("J17").Value = 0.0001038879227879399
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4760280000000001
$ws.Range("N17").Value = 1.428084
$ws.Range("O17").Value = 0.02526391856099382
$ws.Range("P17").Value = 0.02526391856099382
$ws.Range("Q17").Value = 0.001155795984
$ws.Range("R17").Value = 0.010402163856
$ws.Range("S17").Value = 0.000002624616020785329
$ws.Range("T17").Value = 0.000002624616020785329
